# Commit: "updated build_exe.sh and example excel table"
#
# 1) Rename the two worksheets ("kat." -> "cat.")
# 2) On the second sheet ("cat. B"):
#    - bump every data row's height from 18.75pt to 19.5pt (matches the
#      first sheet's row height for its data rows)
#    - shift the "Test 6" record up: it moves from row 5 (where it shared
#      its name in col A with a stray "Test 6" label) so that col A is now
#      blank and the name appears in both B5/C5
#    - fill in two new records in the previously-empty rows 6 and 7:
#      "Test 7" / "Test Name 7" and "Test 8" / "Test Name 8"

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item(1)
$wsB = $wb.Worksheets.Item(2)

$wsA.Name = "cat. A"
$wsB.Name = "cat. B"

# Row heights: every used row on "cat. B" grows from 18.75 to 19.5 points.
for ($r = 1; $r -le 10; $r++) {
    $wsB.Rows.Item($r).RowHeight = 19.5
}

# Row 5: clear col A, keep/refresh the name in B & C.
$wsB.Range("A5").ClearContents()
$wsB.Range("B5").Value = "Test Name 6"
$wsB.Range("C5").Value = "Test Name 6"

# Row 6: new record "Test 7".
$wsB.Range("A6").Value = "Test 7"
$wsB.Range("B6").Value = "Test Name 7"
$wsB.Range("C6").Value = "Test Name 7"

# Row 7: new record "Test 8".
$wsB.Range("A7").Value = "Test 8"
$wsB.Range("B7").Value = "Test Name 8"
$wsB.Range("C7").Value = "Test Name 8"
